$wb = $excel.ActiveWorkbook

# Old / new identifiers used throughout the "Generate Report for Handoff" run.
$oldGuid = "4d845771-9ffb-4acf-80ff-2248905c418e"
$newGuid = "bfbc4c4d-e9d6-437f-8e65-06ea7e7b9a8a"
$oldHash = "bb51db3a43dff848d01d9e9cea737292a49ea77e"
$newHash = "2d2670239bd6316ed6499ad2102ffb486d6e6a78"

$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

$newHandoffDate  = "2016-03-24 08:43:16"
$newZhXlfDate    = "2016-03-24 08:43:12"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Update cell values (shared-string backed) --------------------------

# Source markdown file name, same on every sheet's A2.
$wsOverview.Range("A2").Value = $newMdName
$wsZhCn.Range("A2").Value     = $newMdName
$wsDeDe.Range("A2").Value     = $newMdName

# Overview sheet: Latest Handoff Date (D2).
$wsOverview.Range("D2").Value = $newHandoffDate

# zh-cn sheet: Latest Handoff File (D2) + Latest Handoff Datetime (E2).
$wsZhCn.Range("D2").Value = $newZhXlfName
$wsZhCn.Range("E2").Value = $newZhXlfDate

# de-de sheet: Latest Handoff File (D2) + Latest Handoff Datetime (E2, shares
# the same timestamp text as Overview!D2).
$wsDeDe.Range("D2").Value = $newDeXlfName
$wsDeDe.Range("E2").Value = $newHandoffDate

# --- Update hyperlink display text (kept in sync with the cell text) ----
# NB: iterate with foreach so the existing hyperlink object is mutated in
# place (indexed .Item() access on this host creates a new hyperlink
# instead of editing the existing one).

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMdName
    }
}

foreach ($hl in $wsZhCn.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMdName
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = $newZhXlfName
    }
}

foreach ($hl in $wsDeDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = $newMdName
    } elseif ($addr -eq '$D$2') {
        $hl.TextToDisplay = $newDeXlfName
    }
}
